# api va update form trang thai hoat dong khuyen mai hang tang hang
#
# The sheet used to list 3 gift SKUs (code + name). Replace that sample data
# with a single row reflecting the new "status" tracking used by the
# promo-status form, and add a 3rd column ("Trang thai") to hold it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the old sample data rows (rows 3 and 4) --------------------------
# Delete bottom-up so row numbers don't shift under us.
$ws.Rows.Item(4).Delete() | Out-Null
$ws.Rows.Item(3).Delete() | Out-Null

# --- New data row (entered before the new header, matching editing order) --
$ws.Range("A2").Value2 = 60015505
$ws.Range("C2").Value2 = "x"
$ws.Range("B2").Value2 = "Pathways Reading, Writing 2"

# --- New header cell for the status column --------------------------------
# Copy the existing header style (bold font + border) from A1 onto C1 before
# writing its text so all three headers look the same.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null
$ws.Range("C1").Value2 = "Trạng thái "

# --- Column widths (characters, Calibri 11 default font) -------------------
$ws.Columns.Item(1).ColumnWidth = 11.166666666666666
$ws.Columns.Item(2).ColumnWidth = 32.451822916666664
$ws.Columns.Item(3).ColumnWidth = 13.166666666666666

# --- Page setup: switch to portrait printing --------------------------
$ws.PageSetup.Orientation = 1

# --- Selection, as left by the editor ------------------------------------
$ws.Range("M8").Select() | Out-Null

Write-Host "edit complete"
